$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

# Column C ("Förändrad") holds a date value (serial 45212 = 2023-10-13) for every
# data row (2-70). The automatic update refreshed this "last changed" date to
# serial 45221 (2023-10-22) for all of them.
for ($r = 2; $r -le 70; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    if ($cell.Value2 -eq 45212) {
        $cell.Value = 45221
    }
}
